$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1888.8334
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 2085.6428
$ws.Range("K70").Value = 3600
$ws.Range("L70").Value = 6256.928400000001
$ws.Range("M70").Value = -3330
$ws.Range("N70").Value = -6796.928400000001

$ws.Range("H73").Value = 1888.8334
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 2085.6428
$ws.Range("K73").Value = 3600
$ws.Range("L73").Value = 6256.928400000001
$ws.Range("M73").Value = -2664
$ws.Range("N73").Value = -8128.928400000001

$ws.Range("H74").Value = 3750
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 3750
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360

$ws.Range("H137").Value = 1011.37836
$ws.Range("I137").Value = 914.5909
$ws.Range("J137").Value = 1153.3334
$ws.Range("K137").Value = 2743.7727
$ws.Range("L137").Value = 3460.0002
$ws.Range("M137").Value = -193.7727
$ws.Range("N137").Value = -8560.0002

$ws.Range("H138").Value = 3779.3845
$ws.Range("I138").Value = 1912.2222
$ws.Range("J138").Value = 7980.5
$ws.Range("K138").Value = 5736.6666
$ws.Range("L138").Value = 23941.5
$ws.Range("M138").Value = -596.6665999999996
$ws.Range("N138").Value = -34221.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 547397.5
$ws.Range("I32").Value = 5069.3096
$ws.Range("J32").Value = 2618105
$ws.Range("K32").Value = 5069.3096
$ws.Range("L32").Value = 2618105
$ws.Range("M32").Value = -4782.3096
$ws.Range("N32").Value = -2618679

$ws.Range("H61").Value = 1951.8649
$ws.Range("I61").Value = 2018.4482
$ws.Range("J61").Value = 1710.5
$ws.Range("K61").Value = 2018.4482
$ws.Range("L61").Value = 1710.5
$ws.Range("M61").Value = -1806.4482
$ws.Range("N61").Value = -2134.5

$ws.Range("H74").Value = 757.2963
$ws.Range("I74").Value = 633
$ws.Range("J74").Value = 938.0909
$ws.Range("K74").Value = 633
$ws.Range("L74").Value = 938.0909
$ws.Range("M74").Value = 241
$ws.Range("N74").Value = -2686.0909

$ws.Range("H77").Value = 757.2963
$ws.Range("I77").Value = 633
$ws.Range("J77").Value = 938.0909
$ws.Range("K77").Value = 3165
$ws.Range("L77").Value = 4690.4545
$ws.Range("M77").Value = 1203
$ws.Range("N77").Value = -13426.4545

$ws.Range("H132").Value = 33275.156
$ws.Range("I132").Value = 1640.6875
$ws.Range("J132").Value = 64909.625
$ws.Range("K132").Value = 4922.0625
$ws.Range("L132").Value = 194728.875
$ws.Range("M132").Value = -2392.0625
$ws.Range("N132").Value = -199788.875

$ws.Range("H136").Value = 1951.8649
$ws.Range("I136").Value = 2018.4482
$ws.Range("J136").Value = 1710.5
$ws.Range("K136").Value = 6055.3446
$ws.Range("L136").Value = 5131.5
$ws.Range("M136").Value = -3505.3446
$ws.Range("N136").Value = -10231.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 90
$ws.Range("I22").Value = 90
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 90
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 83

$ws.Range("H26").Value = 21294.2
$ws.Range("I26").Value = 19617.75
$ws.Range("K26").Value = 19617.75
$ws.Range("M26").Value = -19325.75

$ws.Range("H113").Value = 5170
$ws.Range("I113").Value = 5170
$ws.Range("K113").Value = 5170
$ws.Range("M113").Value = -3000

$ws.Range("H134").Value = 2091.4888
$ws.Range("I134").Value = 1978.4634
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 5935.3902
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -3400.3902
$ws.Range("N134").Value = -14820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1313.174
$ws.Range("I31").Value = 1546.3077
$ws.Range("J31").Value = 1010.1
$ws.Range("K31").Value = 1546.3077
$ws.Range("L31").Value = 1010.1
$ws.Range("M31").Value = -1251.3077
$ws.Range("N31").Value = -1600.1

$ws.Range("H34").Value = 1313.174
$ws.Range("I34").Value = 1546.3077
$ws.Range("J34").Value = 1010.1
$ws.Range("K34").Value = 1546.3077
$ws.Range("L34").Value = 1010.1
$ws.Range("M34").Value = -1344.3077
$ws.Range("N34").Value = -1414.1

$ws.Range("H132").Value = 25927974
$ws.Range("I132").Value = 26317572
$ws.Range("J132").Value = 25002676
$ws.Range("K132").Value = 78952716
$ws.Range("L132").Value = 75008028
$ws.Range("M132").Value = -78950186
$ws.Range("N132").Value = -75013088

$ws.Range("H134").Value = 63650110
$ws.Range("I134").Value = 70001120
$ws.Range("J134").Value = 140000
$ws.Range("K134").Value = 210003360
$ws.Range("L134").Value = 420000
$ws.Range("M134").Value = -210000825
$ws.Range("N134").Value = -425070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 777.0599999999999
$ws.Range("I113").Value = 634.7317
$ws.Range("J113").Value = 875.9661
$ws.Range("K113").Value = 1904.1951
$ws.Range("L113").Value = 2627.8983
$ws.Range("M113").Value = 265.8048999999999
$ws.Range("N113").Value = -6967.8983

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1662.72
$ws.Range("I132").Value = 1213.5
$ws.Range("J132").Value = 2461.3333
$ws.Range("K132").Value = 3640.5
$ws.Range("L132").Value = 7383.999899999999
$ws.Range("M132").Value = -1110.5
$ws.Range("N132").Value = -12443.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 5900
$ws.Range("J45").Value = 5900
$ws.Range("L45").Value = 5900
$ws.Range("N45").Value = -6714

$ws.Range("H54").Value = 6380
$ws.Range("J54").Value = 6380
$ws.Range("L54").Value = 6380
$ws.Range("N54").Value = -7668

$ws.Range("H68").Value = 2624.6
$ws.Range("J68").Value = 3083.625
$ws.Range("L68").Value = 3083.625
$ws.Range("N68").Value = -4581.625

$ws.Range("H71").Value = 2624.6
$ws.Range("J71").Value = 3083.625
$ws.Range("L71").Value = 15418.125
$ws.Range("N71").Value = -22906.125

$ws.Range("H132").Value = 9806921
$ws.Range("I132").Value = 11767106
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 35301318
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -35298788
$ws.Range("N132").Value = -23057

$ws.Range("H136").Value = 126775.5
$ws.Range("I136").Value = 126775.5
$ws.Range("K136").Value = 380326.5
$ws.Range("M136").Value = -377776.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0

$ws.Range("M62").ClearContents()
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248

$ws.Range("M65").ClearContents()
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240

$ws.Range("H132").Value = 47620788
$ws.Range("I132").Value = 59830388
$ws.Range("K132").Value = 179491164
$ws.Range("M132").Value = -179488634
